$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    # Force the cell to hold a plain text value (never auto-coerced to a
    # number/float), then restore the cell's original (default) formatting
    # so no stray NumberFormat/Style artifact is left behind.
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "60.615.56"
$ws.Range("E2").Value = "  +1.54%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.606.81"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
Set-TextValue "D5" "575.91"
$ws.Range("E5").Value = "  +2.48%  "

# Row 6 - Solana
Set-TextValue "D6" "143.33"
$ws.Range("E6").Value = "  +0.19%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.25%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.46%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.631.90"
$ws.Range("E9").Value = "  +1.29%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -2.48%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +1.11%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -6.96%  "

# Row 13 - Cardano
Set-TextValue "D13" "0.365"
$ws.Range("E13").Value = "  +2.30%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "3.070.93"
$ws.Range("E14").Value = "  +0.91%  "

# Row 15 - WrappedBTC
Set-TextValue "D15" "60.591.70"
$ws.Range("E15").Value = "  +1.59%  "

# Row 16 - Avalanche
Set-TextValue "D16" "23.28"
$ws.Range("E16").Value = "  -0.05%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +2.60%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.619.41"
$ws.Range("E18").Value = "  +1.13%  "

# Row 19 - Chainlink
Set-TextValue "D19" "11.33"
$ws.Range("E19").Value = "  +8.81%  "

# Row 20 - Polkadot
Set-TextValue "D20" "4.65"
$ws.Range("E20").Value = "  +1.51%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "347.84"
$ws.Range("E21").Value = "  +2.56%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +6.54%  "

# Row 24 - Polygon
Set-TextValue "D24" "0.529"
$ws.Range("E24").Value = "  +11.74%  "

# Row 25 - Litecoin
$ws.Range("E25").Value = "  -0.51%  "

# Row 26 - Binance-PegBSC-USD
Set-TextValue "D26" "0.997"
$ws.Range("E26").Value = "  -0.10%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  -0.86%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue "D28" "7.76"
$ws.Range("E28").Value = "  +3.90%  "

# Row 29 - PEPE
$ws.Range("E29").Value = "  +1.15%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +10.60%  "

# Row 31 - Aptos
Set-TextValue "D31" "6.39"
$ws.Range("E31").Value = "  +2.89%  "

# Row 32 - USDe
$ws.Range("E32").Value = "  -0.11%  "

# Row 33 - Monero
Set-TextValue "D33" "162.42"
$ws.Range("E33").Value = "  +2.86%  "

# Row 34 - EthereumClassic
$ws.Range("E34").Value = "  +2.13%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  +4.10%  "

# Row 36 - Fetch.AI
Set-TextValue "D36" "0.982"
$ws.Range("E36").Value = "  +8.84%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +4.53%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  +7.37%  "

# Row 39 - OKB
$ws.Range("E39").Value = "  +1.24%  "

# Row 40 - Filecoin
$ws.Range("E40").Value = "  +4.23%  "

# Row 41 - SuiNetwork
Set-TextValue "D41" "0.841"
$ws.Range("E41").Value = "  -2.90%  "

# Row 42 - Bittensor
Set-TextValue "D42" "295.93"
$ws.Range("E42").Value = "  +0.15%  "

# Row 43 - Aave
Set-TextValue "D43" "137.54"
$ws.Range("E43").Value = "  -0.36%  "

# Row 44 - FirstDigitalUSD
Set-TextValue "D44" "0.995"
$ws.Range("E44").Value = "  -0.41%  "

# Row 45 - Stellar
Set-TextValue "D45" "0.0986"
$ws.Range("E45").Value = "  +0.81%  "

# Row 46 - Mantle
$ws.Range("E46").Value = "  +2.06%  "

# Row 47 - EnergySwap
Set-TextValue "D47" "19.79"
$ws.Range("E47").Value = "  +3.47%  "

# Row 48 - RenderToken
Set-TextValue "D48" "4.97"
$ws.Range("E48").Value = "  +9.54%  "

# Row 49 - Hedera
$ws.Range("E49").Value = "  +2.11%  "

# Row 51 - WhiteBITCoin -> InjectiveProtocol
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D51" "19.83"
$ws.Range("E51").Value = "  +5.64%  "
